# Updated cryptos list on Mon May 20 20:17:27 UTC 2024 with GitHub Actions
#
# Helper: write a plain-text value into a cell without letting Excel's COM
# layer auto-coerce numeric-looking strings (e.g. "1.00", "2.50") into real
# numbers, which would silently drop formatting like trailing zeros.
# We temporarily force a text number format for the assignment, then reset
# the cell's style back to "Normal" so no extra style survives the edit
# (the source file stores these cells with the default/no style).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "69.722.82"
$ws.Range("E2").Value = "  +5.21%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.430.05"
$ws.Range("E3").Value = "  +11.42%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "586.15"
$ws.Range("E5").Value = "  +1.95%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "183.96"
$ws.Range("E6").Value = "  +8.64%  "

# Row 7 - USDC
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - LidoStakedEther
Set-TextValue $ws.Range("D8") "3.424.95"
$ws.Range("E8").Value = "  +11.36%  "

# Row 9 - XRP
Set-TextValue $ws.Range("D9") "0.531"
$ws.Range("E9").Value = "  +4.26%  "

# Row 10 - Toncoin
Set-TextValue $ws.Range("D10") "6.56"
$ws.Range("E10").Value = "  +3.65%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.156"
$ws.Range("E11").Value = "  +4.79%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("D12") "0.486"
$ws.Range("E12").Value = "  +3.38%  "

# Row 13 - Avalanche
Set-TextValue $ws.Range("D13") "38.21"
$ws.Range("E13").Value = "  +7.12%  "

# Row 14 - ShibaInu
Set-TextValue $ws.Range("D14") "0.0000248"
$ws.Range("E14").Value = "  +3.82%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "4.009.02"
$ws.Range("E15").Value = "  +11.73%  "

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "69.805.03"
$ws.Range("E16").Value = "  +5.50%  "

# Row 17 - TRON (E only)
$ws.Range("E17").Value = "  +1.23%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "3.447.53"
$ws.Range("E18").Value = "  +12.26%  "

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "7.36"
$ws.Range("E19").Value = "  +5.97%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "16.98"
$ws.Range("E20").Value = "  +1.68%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "499.77"
$ws.Range("E21").Value = "  +2.85%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "8.57"
$ws.Range("E22").Value = "  +11.03%  "

# Row 23 - Polygon
Set-TextValue $ws.Range("D23") "0.723"
$ws.Range("E23").Value = "  +5.39%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "86.35"
$ws.Range("E24").Value = "  +4.65%  "

# Row 25 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D25") "13.18"
$ws.Range("E25").Value = "  +4.55%  "

# Row 26 - Fetch.AI
Set-TextValue $ws.Range("D26") "2.37"
$ws.Range("E26").Value = "  +7.81%  "

# Row 27 - RenderToken
Set-TextValue $ws.Range("D27") "10.69"
$ws.Range("E27").Value = "  +4.81%  "

# Row 28 - Dai
Set-TextValue $ws.Range("D28") "0.999"
$ws.Range("E28").Value = "  -0.10%  "

# Row 29 - NEARProtocol
Set-TextValue $ws.Range("D29") "8.22"
$ws.Range("E29").Value = "  +4.83%  "

# Row 30 - ImmutableX
Set-TextValue $ws.Range("D30") "2.50"
$ws.Range("E30").Value = "  +10.88%  "

# Row 31 - PancakeSwap (E only)
$ws.Range("E31").Value = "  +4.08%  "

# Row 32 - EthereumClassic
Set-TextValue $ws.Range("D32") "30.04"
$ws.Range("E32").Value = "  +8.42%  "

# Row 33 - PEPE (E only)
$ws.Range("E33").Value = "  +15.61%  "

# Row 34 - Hedera (E only)
$ws.Range("E34").Value = "  +4.38%  "

# Row 35 - FirstDigitalUSD (E only)
$ws.Range("E35").Value = "  +0.31%  "

# Row 36 - Filecoin
Set-TextValue $ws.Range("D36") "6.07"
$ws.Range("E36").Value = "  +9.07%  "

# Row 37 - Mantle
Set-TextValue $ws.Range("D37") "1.01"
$ws.Range("E37").Value = "  +7.04%  "

# Row 38 - Arweave
Set-TextValue $ws.Range("D38") "48.04"
$ws.Range("E38").Value = "  +2.38%  "

# Row 39 - TheGraph
Set-TextValue $ws.Range("D39") "0.328"
$ws.Range("E39").Value = "  +9.49%  "

# Row 40 - Stacks
Set-TextValue $ws.Range("D40") "2.09"
$ws.Range("E40").Value = "  +6.56%  "

# Row 41 - Kaspa
Set-TextValue $ws.Range("D41") "0.127"
$ws.Range("E41").Value = "  +4.58%  "

# Row 42 - OKB
Set-TextValue $ws.Range("D42") "50.15"
$ws.Range("E42").Value = "  +2.19%  "

# Row 43 - Cosmos
Set-TextValue $ws.Range("D43") "8.65"
$ws.Range("E43").Value = "  +4.68%  "

# Row 44/45 - dogwifhat and Bittensor swap places (with new price/volume data)
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D44") "413.21"
$ws.Range("E44").Value = "  +13.69%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D45") "2.80"
$ws.Range("E45").Value = "  +12.15%  "

# Row 46 - Maker
Set-TextValue $ws.Range("D46") "2.941.48"
$ws.Range("E46").Value = "  +5.70%  "

# Row 47 - InjectiveProtocol
Set-TextValue $ws.Range("D47") "27.86"
$ws.Range("E47").Value = "  +14.55%  "

# Row 48 - VeChain
Set-TextValue $ws.Range("D48") "0.0359"
$ws.Range("E48").Value = "  +4.53%  "

# Row 49 - Monero
Set-TextValue $ws.Range("D49") "134.87"
$ws.Range("E49").Value = "  +0.27%  "

# Row 51 - ThetaToken
Set-TextValue $ws.Range("D51") "2.42"
$ws.Range("E51").Value = "  +12.66%  "
